# "Update country data files" -- Papua New Guinea MSME summary workbook.
#
# 1. Rename the sheet "Data" -> "Summary".
# 2. Remove the disaggregation row (Micro / SMEs / MSMEs headers in B5:D5)
#    that this country's sheet doesn't use, shrinking the used range back
#    down to A1:A3 and pruning the now-unused shared strings.
# 3. Register the "title_" (bold+underline) named cell style alongside the
#    existing name/title/source/HyperLink styles used by this template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "Summary"

# 2. Drop row 5 (B5:D5 = Micro/SMEs/MSMEs) entirely so the sheet's used
#    range goes back to A1:A3 and the shared-string table loses the three
#    now-orphaned entries.
$ws.Rows(5).Delete()

# 3. Add the "title_" named style (bold + underline) to the workbook's
#    style catalog, matching the other name/title/source/HyperLink styles.
$styles = $wb.Styles
$titleStyle = $styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true
